$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows at row 91, shifting existing data (rows 91-163) down to 99-171
$ws.Range("A91:I98").Insert()

# Populate the newly inserted rows with historical data (2019-11-18 .. 2019-11-28)
# Row 91: 2019-11-18
$ws.Cells.Item(91,1).Value = 1574035200
$ws.Cells.Item(91,2).NumberFormat = "@"
$ws.Cells.Item(91,2).Value = "2019-11-18"
$ws.Cells.Item(91,2).NumberFormat = "General"
$ws.Cells.Item(91,3).NumberFormat = "@"
$ws.Cells.Item(91,3).Value = "0210"
$ws.Cells.Item(91,3).NumberFormat = "General"
$ws.Cells.Item(91,4).Value = "KHJB"
$ws.Cells.Item(91,5).Value = 0.29
$ws.Cells.Item(91,6).Value = 0.29
$ws.Cells.Item(91,7).Value = 0.28
$ws.Cells.Item(91,8).Value = 0.28
$ws.Cells.Item(91,9).Value = 61800

# Row 92: 2019-11-19
$ws.Cells.Item(92,1).Value = 1574121600
$ws.Cells.Item(92,2).NumberFormat = "@"
$ws.Cells.Item(92,2).Value = "2019-11-19"
$ws.Cells.Item(92,2).NumberFormat = "General"
$ws.Cells.Item(92,3).NumberFormat = "@"
$ws.Cells.Item(92,3).Value = "0210"
$ws.Cells.Item(92,3).NumberFormat = "General"
$ws.Cells.Item(92,4).Value = "KHJB"
$ws.Cells.Item(92,5).Value = 0.28
$ws.Cells.Item(92,6).Value = 0.28
$ws.Cells.Item(92,7).Value = 0.28
$ws.Cells.Item(92,8).Value = 0.28
$ws.Cells.Item(92,9).Value = 10000

# Row 93: 2019-11-20
$ws.Cells.Item(93,1).Value = 1574208000
$ws.Cells.Item(93,2).NumberFormat = "@"
$ws.Cells.Item(93,2).Value = "2019-11-20"
$ws.Cells.Item(93,2).NumberFormat = "General"
$ws.Cells.Item(93,3).NumberFormat = "@"
$ws.Cells.Item(93,3).Value = "0210"
$ws.Cells.Item(93,3).NumberFormat = "General"
$ws.Cells.Item(93,4).Value = "KHJB"
$ws.Cells.Item(93,5).Value = 0.28
$ws.Cells.Item(93,6).Value = 0.28
$ws.Cells.Item(93,7).Value = 0.275
$ws.Cells.Item(93,8).Value = 0.275
$ws.Cells.Item(93,9).Value = 471300

# Row 94: 2019-11-21
$ws.Cells.Item(94,1).Value = 1574294400
$ws.Cells.Item(94,2).NumberFormat = "@"
$ws.Cells.Item(94,2).Value = "2019-11-21"
$ws.Cells.Item(94,2).NumberFormat = "General"
$ws.Cells.Item(94,3).NumberFormat = "@"
$ws.Cells.Item(94,3).Value = "0210"
$ws.Cells.Item(94,3).NumberFormat = "General"
$ws.Cells.Item(94,4).Value = "KHJB"
$ws.Cells.Item(94,5).Value = 0.28
$ws.Cells.Item(94,6).Value = 0.28
$ws.Cells.Item(94,7).Value = 0.28
$ws.Cells.Item(94,8).Value = 0.28
$ws.Cells.Item(94,9).Value = 245200

# Row 95: 2019-11-22
$ws.Cells.Item(95,1).Value = 1574380800
$ws.Cells.Item(95,2).NumberFormat = "@"
$ws.Cells.Item(95,2).Value = "2019-11-22"
$ws.Cells.Item(95,2).NumberFormat = "General"
$ws.Cells.Item(95,3).NumberFormat = "@"
$ws.Cells.Item(95,3).Value = "0210"
$ws.Cells.Item(95,3).NumberFormat = "General"
$ws.Cells.Item(95,4).Value = "KHJB"
$ws.Cells.Item(95,5).Value = 0.275
$ws.Cells.Item(95,6).Value = 0.275
$ws.Cells.Item(95,7).Value = 0.275
$ws.Cells.Item(95,8).Value = 0.275
$ws.Cells.Item(95,9).Value = 8000

# Row 96: 2019-11-25
$ws.Cells.Item(96,1).Value = 1574640000
$ws.Cells.Item(96,2).NumberFormat = "@"
$ws.Cells.Item(96,2).Value = "2019-11-25"
$ws.Cells.Item(96,2).NumberFormat = "General"
$ws.Cells.Item(96,3).NumberFormat = "@"
$ws.Cells.Item(96,3).Value = "0210"
$ws.Cells.Item(96,3).NumberFormat = "General"
$ws.Cells.Item(96,4).Value = "KHJB"
$ws.Cells.Item(96,5).Value = 0.28
$ws.Cells.Item(96,6).Value = 0.28
$ws.Cells.Item(96,7).Value = 0.275
$ws.Cells.Item(96,8).Value = 0.275
$ws.Cells.Item(96,9).Value = 99700

# Row 97: 2019-11-26
$ws.Cells.Item(97,1).Value = 1574726400
$ws.Cells.Item(97,2).NumberFormat = "@"
$ws.Cells.Item(97,2).Value = "2019-11-26"
$ws.Cells.Item(97,2).NumberFormat = "General"
$ws.Cells.Item(97,3).NumberFormat = "@"
$ws.Cells.Item(97,3).Value = "0210"
$ws.Cells.Item(97,3).NumberFormat = "General"
$ws.Cells.Item(97,4).Value = "KHJB"
$ws.Cells.Item(97,5).Value = 0.275
$ws.Cells.Item(97,6).Value = 0.275
$ws.Cells.Item(97,7).Value = 0.275
$ws.Cells.Item(97,8).Value = 0.275
$ws.Cells.Item(97,9).Value = 112000

# Row 98: 2019-11-28
$ws.Cells.Item(98,1).Value = 1574899200
$ws.Cells.Item(98,2).NumberFormat = "@"
$ws.Cells.Item(98,2).Value = "2019-11-28"
$ws.Cells.Item(98,2).NumberFormat = "General"
$ws.Cells.Item(98,3).NumberFormat = "@"
$ws.Cells.Item(98,3).Value = "0210"
$ws.Cells.Item(98,3).NumberFormat = "General"
$ws.Cells.Item(98,4).Value = "KHJB"
$ws.Cells.Item(98,5).Value = 0.275
$ws.Cells.Item(98,6).Value = 0.285
$ws.Cells.Item(98,7).Value = 0.275
$ws.Cells.Item(98,8).Value = 0.285
$ws.Cells.Item(98,9).Value = 20000
